$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New C2 value (previously empty)
$ws.Range("C2").Value = 7.193183327378438

# Updated values in column E and column C
$ws.Range("E2").Value = 9.591339540850875

$ws.Range("C3").Value = -14.96173956806345
$ws.Range("E3").Value = -4.932343798304595

$ws.Range("C4").Value = 2.682935444832424
$ws.Range("E4").Value = -2.225127715916653

$ws.Range("C5").Value = 9.399485634179229
$ws.Range("E5").Value = 1.811802132286955

$ws.Range("C6").Value = 5.169490031659674
$ws.Range("E6").Value = 9.213376886330305

$ws.Range("C7").Value = -0.3722371047999662
$ws.Range("E7").Value = 2.684220738731935

$ws.Range("C8").Value = 4.098801479368341
$ws.Range("E8").Value = 2.548306621254004

$ws.Range("C9").Value = 3.75051862559701
$ws.Range("E9").Value = 2.714258593289975

$ws.Range("C10").Value = 2.352205130086071
$ws.Range("E10").Value = 3.873414041014778

$ws.Range("C11").Value = 4.083548352538369
$ws.Range("E11").Value = 3.586256146074462

$ws.Range("C12").Value = 4.861590900330692
$ws.Range("E12").Value = 3.297472770389764

$ws.Range("C13").Value = 1.787861866846807
$ws.Range("E13").Value = 4.088367525047842

$ws.Range("C14").Value = -2.21482332957591
$ws.Range("E14").Value = -0.6322362079330346

$ws.Range("C15").Value = 6.09521976277807
$ws.Range("E15").Value = 1.839905110456375

$ws.Range("C16").Value = 3.616930127707629
$ws.Range("E16").Value = 1.391416039405691

$ws.Range("C17").Value = 0.7171092762090492
$ws.Range("E17").Value = 2.755142438739822

$ws.Range("C18").Value = -0.1521036778360019
$ws.Range("E18").Value = 1.645968204809645

$ws.Range("C19").Value = -2.051528019634985
$ws.Range("E19").Value = -0.3224191428759626
